$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated RAD Notice Number Error Message - refresh the "Date" column (B)
# timestamps for the affected rows with newly generated execution times.
$newValues = @(
    @{ Row = 2;  Value = "Wed Dec 20 12:51:39 EST 2023" },
    @{ Row = 3;  Value = "Wed Dec 20 12:51:52 EST 2023" },
    @{ Row = 4;  Value = "Wed Dec 20 12:52:04 EST 2023" },
    @{ Row = 5;  Value = "Wed Dec 20 12:52:16 EST 2023" },
    @{ Row = 6;  Value = "Wed Dec 20 12:52:29 EST 2023" },
    @{ Row = 7;  Value = "Wed Dec 20 12:52:41 EST 2023" },
    @{ Row = 8;  Value = "Wed Dec 20 12:52:53 EST 2023" },
    @{ Row = 9;  Value = "Wed Dec 20 12:53:06 EST 2023" },
    @{ Row = 10; Value = "Wed Dec 20 12:53:18 EST 2023" },
    @{ Row = 11; Value = "Wed Dec 20 12:53:30 EST 2023" },
    @{ Row = 12; Value = "Wed Dec 20 12:53:43 EST 2023" },
    @{ Row = 13; Value = "Wed Dec 20 12:53:55 EST 2023" },
    @{ Row = 14; Value = "Wed Dec 20 12:54:07 EST 2023" },
    @{ Row = 15; Value = "Wed Dec 20 12:54:19 EST 2023" },
    @{ Row = 16; Value = "Wed Dec 20 12:54:32 EST 2023" },
    @{ Row = 17; Value = "Wed Dec 20 12:54:44 EST 2023" },
    @{ Row = 18; Value = "Wed Dec 20 12:54:57 EST 2023" },
    @{ Row = 19; Value = "Wed Dec 20 12:55:09 EST 2023" },
    @{ Row = 20; Value = "Wed Dec 20 12:55:21 EST 2023" },
    @{ Row = 28; Value = "Wed Dec 20 12:55:34 EST 2023" },
    @{ Row = 29; Value = "Wed Dec 20 12:55:47 EST 2023" },
    @{ Row = 30; Value = "Wed Dec 20 12:55:59 EST 2023" },
    @{ Row = 31; Value = "Wed Dec 20 12:56:11 EST 2023" },
    @{ Row = 32; Value = "Wed Dec 20 12:56:24 EST 2023" },
    @{ Row = 33; Value = "Wed Dec 20 12:56:36 EST 2023" }
)

foreach ($entry in $newValues) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Value
}
